# Adds sample data rows to the foresatt, barn and soknad sheets,
# exercising the new "barnehage" (daycare) drop-down list that was
# added to the soknad form (column J / barnehager_prioritert).

$wb = $excel.ActiveWorkbook

# Helper: make a cell hold literal TEXT instead of letting Excel
# auto-convert numeric-looking / date-looking strings into numbers
# or dates. We temporarily force a text number format, assign the
# value, then restore the cell to the default "Normal" style so the
# saved file does not carry a stray custom style on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Helper: give a cell the same look as the bold/boxed header cells
# (style used on column A "index" cells throughout the workbook) by
# copying the format from an existing header cell.
function Copy-HeaderFormat($ws, $headerCell, $targetCell) {
    $headerCell.Copy()
    $targetCell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------
# Sheet "foresatt"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("foresatt")

Copy-HeaderFormat $ws1 $ws1.Cells.Item(1, 2) $ws1.Cells.Item(2, 1)
$ws1.Cells.Item(2, 1).Value = 0
$ws1.Cells.Item(2, 2).Value = 2
$ws1.Cells.Item(2, 3).Value = "a"
$ws1.Cells.Item(2, 4).Value = "a"
$ws1.Cells.Item(2, 5).Value = "s2s"
$ws1.Cells.Item(2, 6).Value = 232

Copy-HeaderFormat $ws1 $ws1.Cells.Item(1, 2) $ws1.Cells.Item(3, 1)
$ws1.Cells.Item(3, 1).Value = 1
$ws1.Cells.Item(3, 2).Value = 1
$ws1.Cells.Item(3, 3).Value = "a"
$ws1.Cells.Item(3, 4).Value = "a"
Set-TextValue $ws1.Cells.Item(3, 5) "1176187"
$ws1.Cells.Item(3, 6).Value = 3242

# ---------------------------------------------------------------
# Sheet "barn"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("barn")

Copy-HeaderFormat $ws3 $ws3.Cells.Item(1, 2) $ws3.Cells.Item(2, 1)
$ws3.Cells.Item(2, 1).Value = 0
$ws3.Cells.Item(2, 2).Value = 1
$ws3.Cells.Item(2, 3).Value = 232323

# ---------------------------------------------------------------
# Sheet "soknad" -- new application row filled in using the
# daycare drop-down (barnehager_prioritert = "Sunshine Preschool")
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("soknad")

Copy-HeaderFormat $ws4 $ws4.Cells.Item(1, 2) $ws4.Cells.Item(2, 1)
$ws4.Cells.Item(2, 1).Value = 0
$ws4.Cells.Item(2, 2).Value = 1
$ws4.Cells.Item(2, 3).Value = 2
$ws4.Cells.Item(2, 4).Value = 2
$ws4.Cells.Item(2, 5).Value = 1
$ws4.Cells.Item(2, 6).Value = "on"
$ws4.Cells.Item(2, 7).Value = ""
$ws4.Cells.Item(2, 8).Value = ""
$ws4.Cells.Item(2, 9).Value = ""
$ws4.Cells.Item(2, 10).Value = "Sunshine Preschool"
$ws4.Cells.Item(2, 11).Value = ""
Set-TextValue $ws4.Cells.Item(2, 12) "2024-10-24"
$ws4.Cells.Item(2, 13).Value = 9000000
